# Auto-generated Excel COM-interop script
# Applies scheduled-runner price/profit updates to the Leve profit sheets
# (ALC, ARM, BSM, CRP, GSM, LTW, WVR) per the commit diff.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$updates = @{
    "H6" = 879.5
    "I6" = 818.6667
    "J6" = 929.2727
    "K6" = 2456.0001
    "L6" = 2787.8181
    "M6" = -2344.0001
    "N6" = -3011.8181
    "H9" = 1332421
    "J9" = 1998463.6
    "L9" = 1998463.6
    "N9" = -1998801.6
    "H62" = 97439280
    "I62" = 115154424
    "K62" = 115154424
    "M62" = -115153800
    "H65" = 97439280
    "I65" = 115154424
    "K65" = 575772120
    "M65" = -575769000
    "H70" = 1941117.6
    "J70" = 1941117.6
    "L70" = 5823352.800000001
    "N70" = -5823892.800000001
    "H73" = 1941117.6
    "J73" = 1941117.6
    "L73" = 5823352.800000001
    "N73" = -5825224.800000001
    "H80" = 1295.2667
    "I80" = 559.3333
    "K80" = 1677.9999
    "M80" = -679.9999
    "H83" = 1295.2667
    "I83" = 559.3333
    "K83" = 5033.9997
    "M83" = -41.9997000000003
    "H133" = 99999.5
    "J133" = 99999.5
    "L133" = 99999.5
    "N133" = -110119.5
    "H137" = 10057.19
    "I137" = 4482.25
    "K137" = 13446.75
    "M137" = -10896.75
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$updates = @{
    "H110" = 12502757
    "I110" = 27779238
    "K110" = 27779238
    "M110" = -27777193
    "H122" = 3157.2559
    "I122" = 2632.2354
    "K122" = 7896.706200000001
    "M122" = -5446.706200000001
    "H132" = 21282438
    "I132" = 25645056
    "J132" = 14676.625
    "K132" = 76935168
    "L132" = 44029.875
    "M132" = -76932638
    "N132" = -49089.875
    "H135" = 151850
    "J135" = 151850
    "L135" = 151850
    "N135" = -161990
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$updates = @{
    "H22" = 2248.1538
    "I22" = 1587.6
    "K22" = 1587.6
    "M22" = -1414.6
    "H107" = 1800.28
    "I107" = 1316.6316
    "K107" = 1316.6316
    "M107" = 603.3684000000001
    "H134" = 19235642
    "J134" = 9942.182000000001
    "L134" = 29826.546
    "N134" = -34896.546
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$updates = @{
    "H16" = 2949.1667
    "I16" = 2523.75
    "K16" = 2523.75
    "M16" = -2236.75
    "H31" = 4572.877
    "I31" = 1928.15
    "K31" = 1928.15
    "M31" = -1633.15
    "H34" = 4572.877
    "I34" = 1928.15
    "K34" = 1928.15
    "M34" = -1726.15
    "H59" = 51856.715
    "I59" = 31249.75
    "J59" = 79332.664
    "K59" = 31249.75
    "L59" = 79332.664
    "M59" = -30104.75
    "N59" = -81622.664
    "H88" = 20666
    "J88" = 20666
    "L88" = 20666
    "N88" = -21478
    "H91" = 20666
    "J91" = 20666
    "L91" = 20666
    "N91" = -23474
    "H99" = 2492.2856
    "I99" = 2530.1538
    "K99" = 2530.1538
    "M99" = -1032.1538
    "H105" = 1296.7858
    "I105" = 1146.2273
    "K105" = 1146.2273
    "M105" = 600.7727
    "H113" = 2949.1667
    "I113" = 2523.75
    "K113" = 2523.75
    "M113" = -353.75
    "H122" = 1925.1364
    "I122" = 1730.2667
    "K122" = 5190.800099999999
    "M122" = -2740.800099999999
    "H126" = 2492.2856
    "I126" = 2530.1538
    "K126" = 7590.4614
    "M126" = -5120.4614
    "H134" = 5219.9375
    "I134" = 4481.9287
    "K134" = 13445.7861
    "M134" = -10910.7861
    "H135" = 450000
    "J135" = 450000
    "L135" = 450000
    "N135" = -460140
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$updates = @{
    "H108" = 149500
    "J108" = 149500
    "L108" = 149500
    "N108" = -157180
    "H132" = 4507.2886
    "I132" = 3646.3618
    "J132" = 12600
    "K132" = 10939.0854
    "L132" = 37800
    "M132" = -8409.0854
    "N132" = -42860
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$updates = @{
    "H111" = 0
    "J111" = 0
    "L111" = 0
    "H132" = 5366.467
    "I132" = 4228.7803
    "J132" = 7821.4736
    "K132" = 12686.3409
    "L132" = 23464.4208
    "M132" = -10156.3409
    "N132" = -28524.4208
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
foreach ($addr in @("N111")) {
    $ws.Range($addr).ClearContents()
}

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$updates = @{
    "H4" = 11500
    "J4" = 11500
    "L4" = 11500
    "N4" = -11726
    "H62" = 0
    "I62" = 0
    "K62" = 0
    "H65" = 0
    "I65" = 0
    "K65" = 0
    "H81" = 12569.293
    "I81" = 6531.8823
    "J81" = 16845.791
    "K81" = 13063.7646
    "L81" = 33691.582
    "M81" = -12002.7646
    "N81" = -35813.582
    "H84" = 12569.293
    "I84" = 6531.8823
    "J84" = 16845.791
    "K84" = 65318.823
    "L84" = 168457.91
    "M84" = -60014.823
    "N84" = -179065.91
    "H132" = 10874334
    "I132" = 16670596
    "J132" = 6341.3125
    "K132" = 50011788
    "L132" = 19023.9375
    "M132" = -50009258
    "N132" = -24083.9375
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
foreach ($addr in @("M62", "M65")) {
    $ws.Range($addr).ClearContents()
}
